$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 56
$ws.Cells.Item(56, 1).Value = 0
$ws.Cells.Item(56, 2).Value = 7.171333983999999
$ws.Cells.Item(56, 3).Value = 1.097692788614881

# New row 57
$ws.Cells.Item(57, 1).Value = 0
$ws.Cells.Item(57, 2).Value = 7.55965918
$ws.Cells.Item(57, 3).Value = 0.9411200094452383
